$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 105
$arr105 = New-Object "object[,]" 1,28
$arr105[0,0] = 6840469
$arr105[0,1] = "Northern Ireland Premier"
$arr105[0,2] = "Northern Ireland Premier"
$arr105[0,3] = 45255.5
$arr105[0,4] = "Linfield"
$arr105[0,5] = "Ballymena Utd"
$arr105[0,6] = 4
$arr105[0,7] = 0
$arr105[0,8] = "H"
$arr105[0,9] = 1.142
$arr105[0,10] = 7.5
$arr105[0,11] = 11
$arr105[0,12] = 1.1
$arr105[0,13] = 8.5
$arr105[0,14] = 21
$arr105[0,15] = -2.5
$arr105[0,16] = 1.875
$arr105[0,17] = 1.925
$arr105[0,18] = 3.5
$arr105[0,19] = 1.95
$arr105[0,20] = 1.85
$arr105[0,21] = 0.1000000000000001
$arr105[0,22] = -1
$arr105[0,23] = -1
$arr105[0,24] = 0.875
$arr105[0,25] = -1
$arr105[0,26] = 0.95
$arr105[0,27] = -1
$ws.Range("B105:AC105").Value = $arr105

# Row 106
$arr106 = New-Object "object[,]" 1,28
$arr106[0,0] = 6840331
$arr106[0,1] = "Northern Ireland Premier"
$arr106[0,2] = "Northern Ireland Premier"
$arr106[0,3] = 45255.5
$arr106[0,4] = "Coleraine"
$arr106[0,5] = "Carrick Rangers"
$arr106[0,6] = 1
$arr106[0,7] = 1
$arr106[0,8] = "D"
$arr106[0,9] = 1.5
$arr106[0,10] = 4
$arr106[0,11] = 5
$arr106[0,12] = 1.6
$arr106[0,13] = 4
$arr106[0,14] = 4.5
$arr106[0,15] = -1
$arr106[0,16] = 2
$arr106[0,17] = 1.8
$arr106[0,18] = 2.75
$arr106[0,19] = 1.825
$arr106[0,20] = 1.975
$arr106[0,21] = -1
$arr106[0,22] = 3
$arr106[0,23] = -1
$arr106[0,24] = -1
$arr106[0,25] = 0.8
$arr106[0,26] = -1
$arr106[0,27] = 0.9750000000000001
$ws.Range("B106:AC106").Value = $arr106

# Row 107
$arr107 = New-Object "object[,]" 1,28
$arr107[0,0] = 6840467
$arr107[0,1] = "Northern Ireland Premier"
$arr107[0,2] = "Northern Ireland Premier"
$arr107[0,3] = 45255.5
$arr107[0,4] = "Cliftonville"
$arr107[0,5] = "Crusaders"
$arr107[0,6] = 3
$arr107[0,7] = 0
$arr107[0,8] = "H"
$arr107[0,9] = 2.25
$arr107[0,10] = 3.25
$arr107[0,11] = 2.75
$arr107[0,12] = 2.1
$arr107[0,13] = 3.25
$arr107[0,14] = 3.3
$arr107[0,15] = -0.25
$arr107[0,16] = 1.825
$arr107[0,17] = 1.975
$arr107[0,18] = 2.5
$arr107[0,19] = 1.875
$arr107[0,20] = 1.925
$arr107[0,21] = 1.1
$arr107[0,22] = -1
$arr107[0,23] = -1
$arr107[0,24] = 0.825
$arr107[0,25] = -1
$arr107[0,26] = 0.875
$arr107[0,27] = -1
$ws.Range("B107:AC107").Value = $arr107

# Row 121
$arr121 = New-Object "object[,]" 1,28
$arr121[0,0] = 6840475
$arr121[0,1] = "Northern Ireland Premier"
$arr121[0,2] = "Northern Ireland Premier"
$arr121[0,3] = 45275.69791666666
$arr121[0,4] = "Ballymena Utd"
$arr121[0,5] = "Carrick Rangers"
$arr121[0,6] = 0
$arr121[0,7] = 2
$arr121[0,8] = "A"
$arr121[0,9] = 3.4
$arr121[0,10] = 3.3
$arr121[0,11] = 2
$arr121[0,12] = 3.5
$arr121[0,13] = 3.25
$arr121[0,14] = 2
$arr121[0,15] = 0.5
$arr121[0,16] = 1.775
$arr121[0,17] = 2.025
$arr121[0,18] = 2.5
$arr121[0,19] = 2.025
$arr121[0,20] = 1.775
$arr121[0,21] = -1
$arr121[0,22] = -1
$arr121[0,23] = 1
$arr121[0,24] = -1
$arr121[0,25] = 1.025
$arr121[0,26] = -1
$arr121[0,27] = 0.7749999999999999
$ws.Range("B121:AC121").Value = $arr121

# Row 122
$arr122 = New-Object "object[,]" 1,28
$arr122[0,0] = 6840337
$arr122[0,1] = "Northern Ireland Premier"
$arr122[0,2] = "Northern Ireland Premier"
$arr122[0,3] = 45275.69791666666
$arr122[0,4] = "Crusaders"
$arr122[0,5] = "Larne FC"
$arr122[0,6] = 0
$arr122[0,7] = 3
$arr122[0,8] = "A"
$arr122[0,9] = 3.1
$arr122[0,10] = 3.2
$arr122[0,11] = 2.25
$arr122[0,12] = 3.5
$arr122[0,13] = 3.2
$arr122[0,14] = 2.05
$arr122[0,15] = 0.25
$arr122[0,16] = 2
$arr122[0,17] = 1.8
$arr122[0,18] = 2.25
$arr122[0,19] = 2.025
$arr122[0,20] = 1.775
$arr122[0,21] = -1
$arr122[0,22] = -1
$arr122[0,23] = 1.05
$arr122[0,24] = -1
$arr122[0,25] = 0.8
$arr122[0,26] = 1.025
$arr122[0,27] = -1
$ws.Range("B122:AC122").Value = $arr122

# Row 162
$arr162 = New-Object "object[,]" 1,28
$arr162[0,0] = 6840297
$arr162[0,1] = "Northern Ireland Premier"
$arr162[0,2] = "Northern Ireland Premier"
$arr162[0,3] = 45331.69791666666
$arr162[0,4] = "Newry City"
$arr162[0,5] = "Ballymena Utd"
$arr162[0,6] = 1
$arr162[0,7] = 1
$arr162[0,8] = "D"
$arr162[0,9] = 2.625
$arr162[0,10] = 3.3
$arr162[0,11] = 2.45
$arr162[0,12] = 2.75
$arr162[0,13] = 3.3
$arr162[0,14] = 2.375
$arr162[0,15] = 0
$arr162[0,16] = 2
$arr162[0,17] = 1.8
$arr162[0,18] = 2.25
$arr162[0,19] = 1.775
$arr162[0,20] = 2.025
$arr162[0,21] = -1
$arr162[0,22] = 2.3
$arr162[0,23] = -1
$arr162[0,24] = 0
$arr162[0,25] = -0
$arr162[0,26] = -0.5
$arr162[0,27] = 0.5125
$ws.Range("B162:AC162").Value = $arr162

# Row 163
$arr163 = New-Object "object[,]" 1,28
$arr163[0,0] = 6839226
$arr163[0,1] = "Northern Ireland Premier"
$arr163[0,2] = "Northern Ireland Premier"
$arr163[0,3] = 45331.69791666666
$arr163[0,4] = "Larne FC"
$arr163[0,5] = "Loughgall"
$arr163[0,6] = 2
$arr163[0,7] = 0
$arr163[0,8] = "H"
$arr163[0,9] = 1.181
$arr163[0,10] = 6.5
$arr163[0,11] = 11
$arr163[0,12] = 1.181
$arr163[0,13] = 6.5
$arr163[0,14] = 11
$arr163[0,15] = -2
$arr163[0,16] = 1.925
$arr163[0,17] = 1.875
$arr163[0,18] = 3
$arr163[0,19] = 1.825
$arr163[0,20] = 1.975
$arr163[0,21] = 0.181
$arr163[0,22] = -1
$arr163[0,23] = -1
$arr163[0,24] = 0
$arr163[0,25] = -0
$arr163[0,26] = -1
$arr163[0,27] = 0.9750000000000001
$ws.Range("B163:AC163").Value = $arr163

# Row 170
$arr170 = New-Object "object[,]" 1,28
$arr170[0,0] = 6839225
$arr170[0,1] = "Northern Ireland Premier"
$arr170[0,2] = "Northern Ireland Premier"
$arr170[0,3] = 45339.5
$arr170[0,4] = "Loughgall"
$arr170[0,5] = "Crusaders"
$arr170[0,6] = 0
$arr170[0,7] = 1
$arr170[0,8] = "A"
$arr170[0,9] = 5.25
$arr170[0,10] = 4
$arr170[0,11] = 1.5
$arr170[0,12] = 5
$arr170[0,13] = 3.8
$arr170[0,14] = 1.55
$arr170[0,15] = 1
$arr170[0,16] = 1.8
$arr170[0,17] = 2
$arr170[0,18] = 2.75
$arr170[0,19] = 1.875
$arr170[0,20] = 1.925
$arr170[0,21] = -1
$arr170[0,22] = -1
$arr170[0,23] = 0.55
$arr170[0,24] = 0
$arr170[0,25] = -0
$arr170[0,26] = -1
$arr170[0,27] = 0.925
$ws.Range("B170:AC170").Value = $arr170

# Row 171
$arr171 = New-Object "object[,]" 1,28
$arr171[0,0] = 6840298
$arr171[0,1] = "Northern Ireland Premier"
$arr171[0,2] = "Northern Ireland Premier"
$arr171[0,3] = 45339.5
$arr171[0,4] = "Ballymena Utd"
$arr171[0,5] = "Larne FC"
$arr171[0,6] = 0
$arr171[0,7] = 1
$arr171[0,8] = "A"
$arr171[0,9] = 9.5
$arr171[0,10] = 5
$arr171[0,11] = 1.285
$arr171[0,12] = 7
$arr171[0,13] = 4.2
$arr171[0,14] = 1.4
$arr171[0,15] = 1.25
$arr171[0,16] = 1.85
$arr171[0,17] = 1.95
$arr171[0,18] = 2.5
$arr171[0,19] = 1.95
$arr171[0,20] = 1.85
$arr171[0,21] = -1
$arr171[0,22] = -1
$arr171[0,23] = 0.3999999999999999
$arr171[0,24] = 0.425
$arr171[0,25] = -0.5
$arr171[0,26] = -1
$arr171[0,27] = 0.8500000000000001
$ws.Range("B171:AC171").Value = $arr171

# Row 190
$arr190 = New-Object "object[,]" 1,28
$arr190[0,0] = 6840961
$arr190[0,1] = "Northern Ireland Premier"
$arr190[0,2] = "Northern Ireland Premier"
$arr190[0,3] = 45367.5
$arr190[0,4] = "Loughgall"
$arr190[0,5] = "Cliftonville"
$arr190[0,6] = 2
$arr190[0,7] = 3
$arr190[0,8] = "A"
$arr190[0,9] = 7
$arr190[0,10] = 4.5
$arr190[0,11] = 1.4
$arr190[0,12] = 4.2
$arr190[0,13] = 4
$arr190[0,14] = 1.65
$arr190[0,15] = 0.75
$arr190[0,16] = 1.95
$arr190[0,17] = 1.85
$arr190[0,18] = 3
$arr190[0,19] = 1.925
$arr190[0,20] = 1.875
$arr190[0,21] = -1
$arr190[0,22] = -1
$arr190[0,23] = 0.6499999999999999
$arr190[0,24] = -0.5
$arr190[0,25] = 0.425
$arr190[0,26] = 0.925
$arr190[0,27] = -1
$ws.Range("B190:AC190").Value = $arr190

# Row 191
$arr191 = New-Object "object[,]" 1,28
$arr191[0,0] = 6841445
$arr191[0,1] = "Northern Ireland Premier"
$arr191[0,2] = "Northern Ireland Premier"
$arr191[0,3] = 45367.5
$arr191[0,4] = "Crusaders"
$arr191[0,5] = "Coleraine"
$arr191[0,6] = 1
$arr191[0,7] = 1
$arr191[0,8] = "D"
$arr191[0,9] = 1.75
$arr191[0,10] = 3.75
$arr191[0,11] = 4.2
$arr191[0,12] = 1.8
$arr191[0,13] = 3.6
$arr191[0,14] = 4.2
$arr191[0,15] = -0.75
$arr191[0,16] = 2
$arr191[0,17] = 1.8
$arr191[0,18] = 2.5
$arr191[0,19] = 1.825
$arr191[0,20] = 1.975
$arr191[0,21] = -1
$arr191[0,22] = 2.6
$arr191[0,23] = -1
$arr191[0,24] = -1
$arr191[0,25] = 0.8
$arr191[0,26] = -1
$arr191[0,27] = 0.9750000000000001
$ws.Range("B191:AC191").Value = $arr191

# Row 192
$arr192 = New-Object "object[,]" 1,28
$arr192[0,0] = 6840958
$arr192[0,1] = "Northern Ireland Premier"
$arr192[0,2] = "Northern Ireland Premier"
$arr192[0,3] = 45367.5
$arr192[0,4] = "Ballymena Utd"
$arr192[0,5] = "Carrick Rangers"
$arr192[0,6] = 0
$arr192[0,7] = 2
$arr192[0,8] = "A"
$arr192[0,9] = 3
$arr192[0,10] = 3.4
$arr192[0,11] = 2.2
$arr192[0,12] = 3.6
$arr192[0,13] = 3.5
$arr192[0,14] = 1.95
$arr192[0,15] = 0.5
$arr192[0,16] = 1.825
$arr192[0,17] = 1.975
$arr192[0,18] = 2.5
$arr192[0,19] = 1.975
$arr192[0,20] = 1.825
$arr192[0,21] = -1
$arr192[0,22] = -1
$arr192[0,23] = 0.95
$arr192[0,24] = -1
$arr192[0,25] = 0.9750000000000001
$arr192[0,26] = -1
$arr192[0,27] = 0.825
$ws.Range("B192:AC192").Value = $arr192

Write-Output "done"